$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header strings: "_old" -> "_FV2410" (cols A-J), "_new" -> "_FV2504" (cols L-U) ---
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $oldVal = $cell.Value2
    $cell.Value = ($oldVal -replace '_old$', '_FV2410')
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $oldVal = $cell.Value2
    $cell.Value = ($oldVal -replace '_new$', '_FV2504')
}

# --- Freeze the header row (row 1) ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
[void]($excel.ActiveWindow.FreezePanes = $true)

# --- Turn the data range into an Excel Table ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U56"), $null, 1)
$tbl.Name = "Table1"
